$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the summary totals -----------------------------------------
# VALOR MORA total (E11) grew because a new overdue period (2509) was added
$ws.Range("E11").Value = 2688000
# Cant. Periodos (F13) increases from 19 to 20 periods
$ws.Range("F13").Value = 20

# --- 2. Add a new detail row for period 2509 -------------------------------
# Insert a blank row right after the current last data row (34) so the old
# last row (period 2508) becomes a regular middle row and the brand new row
# (period 2509) becomes the new last row of the table.
$ws.Rows.Item(35).Insert()

# Give the new last row (35) the special "last row" look that the old last
# row (34) used to have.
$ws.Range("B34:J34").Copy()
$ws.Range("B35:J35").PasteSpecial(-4122) | Out-Null

# Give the row that used to be last (34, period 2508) the regular "middle
# row" look shared by all the other detail rows (copy format from row 33).
$ws.Range("B33:J33").Copy()
$ws.Range("B34:J34").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the values for the new period-2509 row.
$ws.Range("B35").Value = "CC"
$ws.Range("C35").Value = "11511691"
$ws.Range("D35").Value = "JAIRO ALFREDO JIMENEZ BARON"
$ws.Range("E35").Value = "2509"
$ws.Range("F35").Value = 134400
$ws.Range("G35").Value = 3360000

# --- 3. Center the "Periodo Mora" column for every detail row -------------
$ws.Range("E16:E35").HorizontalAlignment = -4108
